$d = $word.ActiveDocument

# 1. Add the portfolio link "www.derekhan.info" into the empty paragraph
#    that sits right after the contact-info block (phone number line)
#    and before the "OBJECTIVE STATEMENT" heading.
$p = $d.Paragraphs.Item(4)
if ($p.Range.Text -eq "" -or $p.Range.Text -eq "`r") {
    $r = $p.Range
    $r.Text = "www.derekhan.info"
    $r.Font.Name = "Georgia"
    $r.Font.NameFarEast = "SimSun"
    $r.Font.NameBi = "Times New Roman"
    $r.Font.Size = 10
    $r.Font.SizeBi = 10
    $r.LanguageIDFarEast = "zh-CN"
}

# 2. Click right after "ability" (before " to adapt ...") in the objective
#    statement paragraph. Word records the last edit point with its
#    "_GoBack" bookmark, so this both splits the run at that point and
#    relocates the bookmark from its old spot after "cumulative" to here.
$find = $d.Content
$found = $find.Find.Execute("ability", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $d.Range($find.End, $find.End)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}
